$d = $word.ActiveDocument

# Q1: merge "support " + "iteration?" into one run (no-op text replace forces run merge)
$d.Content.Find.Execute("Q1. Which two operator overloading methods can you use in your classes to support iteration?", $true, $false, $false, $false, $false, $true, 1, $false, "Q1. Which two operator overloading methods can you use in your classes to support iteration?", 2) | Out-Null

# Q2: merge "...printing" + "?" into one run
$d.Content.Find.Execute("Q2. In what contexts do the two operator overloading methods manage printing?", $true, $false, $false, $false, $false, $true, 1, $false, "Q2. In what contexts do the two operator overloading methods manage printing?", 2) | Out-Null

# Q4: merge "...in-place " + "addition?" into one run
$d.Content.Find.Execute("Q4. In a class, how do you capture in-place addition?", $true, $false, $false, $false, $false, $true, 1, $false, "Q4. In a class, how do you capture in-place addition?", 2) | Out-Null

# Answer to Q4: merge "__` method" + " " + "This method is called when you use the `+=` operator" into one run
$q4ans = "__`` method This method is called when you use the ``+=`` operator"
$d.Content.Find.Execute($q4ans, $true, $false, $false, $false, $false, $true, 1, $false, $q4ans, 2) | Out-Null

# Q5: merge "...overloading" + "?" into one run
$d.Content.Find.Execute("Q5. When is it appropriate to use operator overloading?", $true, $false, $false, $false, $false, $true, 1, $false, "Q5. When is it appropriate to use operator overloading?", 2) | Out-Null
